$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A30: empty text (quote-prefixed empty string keeps it a Text cell, not a truly blank one)
$ws.Range("A30").Value = "'"
$ws.Range("A30").Style = "Normal"

$ws.Range("B30").Value = "احمد"

# C30: digit-only text must be forced to Text so it isn't coerced to a number
$ws.Range("C30").Value = "'22"
$ws.Range("C30").Style = "Normal"

$ws.Range("D30").Value = "الصمود"
$ws.Range("E30").Value = "الرحلة 1"
$ws.Range("F30").Value = "C1"
$ws.Range("G30").Value = "WCK"
$ws.Range("H30").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٠٨:١٩ م"
